# Actualización automática 2025-12-01 08:30:07
#
# Applies the monthly "roll forward" update:
#  - Sheet "VENTAS POR GRUPO": a handful of cells are zeroed out (and their
#    related progress counters updated) because the corresponding sale
#    moved out of the tracked window.
#  - Sheet "VENTA MENSUAL": the tracked month columns (C..F) shift one
#    month forward (headers agosto/sept/oct/nov -> sept/oct/nov/dic) and
#    the per-row monthly figures shift left by one column, with the new
#    right-most month (F) reset to 0. Column widths for C/D/E follow the
#    new layout as well.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M3").Value = 0
$wsGrupo.Range("L19").Value = 0
$wsGrupo.Range("L21").Value = 0
$wsGrupo.Range("C24").Value = 0
$wsGrupo.Range("I24").Value = 0
$wsGrupo.Range("L24").Value = 0
$wsGrupo.Range("N24").Value = 0

$wsGrupo.Range("C32").Value = "0 de 30"
$wsGrupo.Range("I32").Value = "0 de 30"
$wsGrupo.Range("L32").Value = "0 de 30"
$wsGrupo.Range("M32").Value = "0 de 30"
$wsGrupo.Range("N32").Value = "0 de 30"

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Month headers roll forward by one month.
$wsMensual.Range("C1").Value = "septiembre"
$wsMensual.Range("D1").Value = "octubre"
$wsMensual.Range("E1").Value = "noviembre"
$wsMensual.Range("F1").Value = "diciembre"

# Column widths follow the same C/D/E layout used for the new months.
$wsMensual.Columns.Item(3).ColumnWidth = 15.17
$wsMensual.Columns.Item(4).ColumnWidth = 12.17
$wsMensual.Columns.Item(5).ColumnWidth = 14.17

# Row data shifts one column to the left (new F becomes 0).
$rows = 2..32
foreach ($r in $rows) {
    $c = $wsMensual.Range("C$r").Value2
    $d = $wsMensual.Range("D$r").Value2
    $e = $wsMensual.Range("E$r").Value2
    $f = $wsMensual.Range("F$r").Value2

    $wsMensual.Range("C$r").Value = $d
    $wsMensual.Range("D$r").Value = $e
    $wsMensual.Range("E$r").Value = $f
    $wsMensual.Range("F$r").Value = 0
}

Write-Output "done"
